# Add manager name to ins_request:
# Fix typo in the "requested risks" label and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tech_info")

# Correct the misspelled label "Заправшиваемые риски:" -> "Запрашиваемые риски:"
$ws.Range("A17").Value = "Запрашиваемые риски:"

# Update the sheet's active cell/selection to A18
$ws.Range("A18").Select()
